$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same underlying data table
# (mirrored sheets). Update the "想去人数" (column F) counters on each.

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F6").Value = 1109
    $ws.Range("F9").Value = 7632
    $ws.Range("F11").Value = 427
    $ws.Range("F12").Value = 353
    $ws.Range("F16").Value = 7779
    $ws.Range("F29").Value = 415
    $ws.Range("F36").Value = 76
}

# F30 differs slightly between the two mirrored sheets
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F30").Value = 857

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F30").Value = 858
